# Auto-generated edit script: updates market-derived profit columns (H-N)
# on several worksheets per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1050.1428
$ws.Range("I32").Value = 915
$ws.Range("J32").Value = 1104.2
$ws.Range("K32").Value = 915
$ws.Range("L32").Value = 1104.2
$ws.Range("M32").Value = -589
$ws.Range("N32").Value = -1756.2
$ws.Range("H41").Value = 1329.8572
$ws.Range("I41").Value = 2899.75
$ws.Range("J41").Value = 701.9
$ws.Range("K41").Value = 2899.75
$ws.Range("L41").Value = 701.9
$ws.Range("M41").Value = -2459.75
$ws.Range("N41").Value = -1581.9
$ws.Range("H129").Value = 4598.8213
$ws.Range("J129").Value = 1067.9615
$ws.Range("L129").Value = 3203.8845
$ws.Range("N129").Value = -13203.8845
$ws.Range("H132").Value = 8069803.5
$ws.Range("I132").Value = 9621208
$ws.Range("J132").Value = 2497.2
$ws.Range("K132").Value = 28863624
$ws.Range("L132").Value = 7491.599999999999
$ws.Range("M132").Value = -28861094
$ws.Range("N132").Value = -12551.6
$ws.Range("H137").Value = 1462.3438
$ws.Range("I137").Value = 1405.68
$ws.Range("K137").Value = 4217.04
$ws.Range("M137").Value = -1667.04
$ws.Range("H141").Value = 5019
$ws.Range("I141").Value = 6047.5
$ws.Range("K141").Value = 18142.5
$ws.Range("M141").Value = -12962.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31712.803
$ws.Range("I32").Value = 4698.06
$ws.Range("J32").Value = 256835.67
$ws.Range("K32").Value = 4698.06
$ws.Range("L32").Value = 256835.67
$ws.Range("M32").Value = -4411.06
$ws.Range("N32").Value = -257409.67
$ws.Range("H61").Value = 1308.0385
$ws.Range("I61").Value = 1191.3182
$ws.Range("J61").Value = 1950
$ws.Range("K61").Value = 1191.3182
$ws.Range("L61").Value = 1950
$ws.Range("M61").Value = -979.3181999999999
$ws.Range("N61").Value = -2374
$ws.Range("H74").Value = 2427
$ws.Range("I74").Value = 1302.4736
$ws.Range("J74").Value = 4801
$ws.Range("K74").Value = 1302.4736
$ws.Range("L74").Value = 4801
$ws.Range("M74").Value = -428.4736
$ws.Range("N74").Value = -6549
$ws.Range("H77").Value = 2427
$ws.Range("I77").Value = 1302.4736
$ws.Range("J77").Value = 4801
$ws.Range("K77").Value = 6512.368
$ws.Range("L77").Value = 24005
$ws.Range("M77").Value = -2144.368
$ws.Range("N77").Value = -32741
$ws.Range("H122").Value = 2697.923
$ws.Range("I122").Value = 2247.75
$ws.Range("J122").Value = 2898
$ws.Range("K122").Value = 6743.25
$ws.Range("L122").Value = 8694
$ws.Range("M122").Value = -4293.25
$ws.Range("N122").Value = -13594
$ws.Range("H132").Value = 1989.24
$ws.Range("I132").Value = 1587.25
$ws.Range("K132").Value = 4761.75
$ws.Range("M132").Value = -2231.75
$ws.Range("H136").Value = 1308.0385
$ws.Range("I136").Value = 1191.3182
$ws.Range("J136").Value = 1950
$ws.Range("K136").Value = 3573.9546
$ws.Range("L136").Value = 5850
$ws.Range("M136").Value = -1023.9546
$ws.Range("N136").Value = -10950

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66850.766
$ws.Range("I86").Value = 111579.6
$ws.Range("J86").Value = 2952.4285
$ws.Range("K86").Value = 111579.6
$ws.Range("L86").Value = 2952.4285
$ws.Range("M86").Value = -110456.6
$ws.Range("N86").Value = -5198.4285
$ws.Range("H89").Value = 66850.766
$ws.Range("I89").Value = 111579.6
$ws.Range("J89").Value = 2952.4285
$ws.Range("K89").Value = 557898
$ws.Range("L89").Value = 14762.1425
$ws.Range("M89").Value = -552282
$ws.Range("N89").Value = -25994.1425
$ws.Range("H134").Value = 2755.5557
$ws.Range("I134").Value = 3018.9
$ws.Range("J134").Value = 2003.1428
$ws.Range("K134").Value = 9056.700000000001
$ws.Range("L134").Value = 6009.428400000001
$ws.Range("M134").Value = -6521.700000000001
$ws.Range("N134").Value = -11079.4284

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22665.592
$ws.Range("I31").Value = 37063.785
$ws.Range("K31").Value = 37063.785
$ws.Range("M31").Value = -36768.785
$ws.Range("H34").Value = 22665.592
$ws.Range("I34").Value = 37063.785
$ws.Range("K34").Value = 37063.785
$ws.Range("M34").Value = -36861.785
$ws.Range("H58").Value = 15182.895
$ws.Range("I58").Value = 1876
$ws.Range("J58").Value = 37994.715
$ws.Range("K58").Value = 1876
$ws.Range("L58").Value = 37994.715
$ws.Range("M58").Value = -1673
$ws.Range("N58").Value = -38400.715
$ws.Range("H122").Value = 598.6667
$ws.Range("I122").Value = 562.8570999999999
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 1688.5713
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = 761.4287000000002
$ws.Range("N122").Value = -8200
$ws.Range("H132").Value = 33336514
$ws.Range("I132").Value = 31252962
$ws.Range("J132").Value = 38465256
$ws.Range("K132").Value = 93758886
$ws.Range("L132").Value = 115395768
$ws.Range("M132").Value = -93756356
$ws.Range("N132").Value = -115400828
$ws.Range("H134").Value = 1347.1111
$ws.Range("I134").Value = 1326.8462
$ws.Range("J134").Value = 1399.8
$ws.Range("K134").Value = 3980.5386
$ws.Range("L134").Value = 4199.4
$ws.Range("M134").Value = -1445.5386
$ws.Range("N134").Value = -9269.4
$ws.Range("H136").Value = 15182.895
$ws.Range("I136").Value = 1876
$ws.Range("J136").Value = 37994.715
$ws.Range("K136").Value = 5628
$ws.Range("L136").Value = 113984.145
$ws.Range("M136").Value = -3078
$ws.Range("N136").Value = -119084.145

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9590.75
$ws.Range("I5").Value = 1405.6666
$ws.Range("J5").Value = 12319.111
$ws.Range("K5").Value = 4216.9998
$ws.Range("L5").Value = 36957.333
$ws.Range("M5").Value = -4104.9998
$ws.Range("N5").Value = -37181.333
$ws.Range("H122").Value = 4370.5386
$ws.Range("I122").Value = 377.9524
$ws.Range("K122").Value = 3401.5716
$ws.Range("M122").Value = -951.5716000000002
$ws.Range("H132").Value = 3513.7896
$ws.Range("J132").Value = 3450.818
$ws.Range("L132").Value = 31057.362
$ws.Range("N132").Value = -36117.362
$ws.Range("H135").Value = 9590.75
$ws.Range("I135").Value = 1405.6666
$ws.Range("J135").Value = 12319.111
$ws.Range("K135").Value = 12650.9994
$ws.Range("L135").Value = 110871.999
$ws.Range("M135").Value = -10115.9994
$ws.Range("N135").Value = -115941.999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1152.8182
$ws.Range("J122").Value = 1354.8
$ws.Range("L122").Value = 4064.4
$ws.Range("N122").Value = -8964.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 11200
$ws.Range("I24").Value = 2000
$ws.Range("J24").Value = 13500
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 13500
$ws.Range("M24").Value = -1657
$ws.Range("N24").Value = -14186
$ws.Range("H132").Value = 1576.3438
$ws.Range("I132").Value = 1301.36
$ws.Range("J132").Value = 2558.4285
$ws.Range("K132").Value = 3904.08
$ws.Range("L132").Value = 7675.2855
$ws.Range("M132").Value = -1374.08
$ws.Range("N132").Value = -12735.2855
$ws.Range("H136").Value = 2238.125
$ws.Range("I136").Value = 2400
$ws.Range("J136").Value = 2076.25
$ws.Range("K136").Value = 7200
$ws.Range("L136").Value = 6228.75
$ws.Range("M136").Value = -4650
$ws.Range("N136").Value = -11328.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H132").Value = 3687.1904
$ws.Range("I132").Value = 3964.8667
$ws.Range("J132").Value = 2993
$ws.Range("K132").Value = 11894.6001
$ws.Range("L132").Value = 8979
$ws.Range("M132").Value = -9364.6001
$ws.Range("N132").Value = -14039
$ws.Range("H136").Value = 1189
$ws.Range("I136").Value = 952
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 2856
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -306

